$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for Espinaca (La Araucania), dated
# 2021-12-14 (serial 44544). It belongs right before the current row 38,
# so insert a fresh row there which pushes every following record down by
# one (old row 38 -> 39, ... old row 111 -> 112).
$ws.Rows.Item(38).Insert()

$ws.Range("A38").Value = 10
$ws.Range("B38").Value = "Vega Modelo de Temuco"
$ws.Range("C38").Value = "La Araucanía"
$ws.Range("D38").Value = 44544
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = 100112012
$ws.Range("G38").Value = "Espinaca"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 65
$ws.Range("K38").Value = 9000
$ws.Range("L38").Value = 9000
$ws.Range("M38").Value = 9000
$ws.Range("N38").Value = "$/docena de atados"
$ws.Range("O38").Value = "Región de La Araucanía"
$ws.Range("P38").Value = 3000
$ws.Range("Q38").Value = 3
$ws.Range("R38").Value = "Hortaliza"
